$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(2, 8).Value = 504
$ws.Cells.Item(2, 9).Value = 488.83334
$ws.Cells.Item(2, 11).Value = 488.83334
$ws.Cells.Item(2, 13).Value = -375.83334
$ws.Cells.Item(19, 8).Value = 1962.4615
$ws.Cells.Item(19, 10).Value = 2414.7
$ws.Cells.Item(19, 12).Value = 2414.7
$ws.Cells.Item(19, 14).Value = -2764.7
$ws.Cells.Item(33, 8).Value = 12921.685
$ws.Cells.Item(33, 9).Value = 14609.5
$ws.Cells.Item(33, 10).Value = 3920
$ws.Cells.Item(33, 11).Value = 14609.5
$ws.Cells.Item(33, 12).Value = 3920
$ws.Cells.Item(33, 13).Value = -14380.5
$ws.Cells.Item(33, 14).Value = -4378
$ws.Cells.Item(64, 8).Value = 8747.1
$ws.Cells.Item(64, 9).Value = 8830.111000000001
$ws.Cells.Item(64, 11).Value = 8830.111000000001
$ws.Cells.Item(64, 13).Value = -8582.111000000001
$ws.Cells.Item(67, 8).Value = 8747.1
$ws.Cells.Item(67, 9).Value = 8830.111000000001
$ws.Cells.Item(67, 11).Value = 8830.111000000001
$ws.Cells.Item(67, 13).Value = -7972.111000000001
$ws.Cells.Item(100, 8).Value = 2350.3572
$ws.Cells.Item(100, 9).Value = 1800.6666
$ws.Cells.Item(100, 10).Value = 3339.8
$ws.Cells.Item(100, 11).Value = 1800.6666
$ws.Cells.Item(100, 12).Value = 3339.8
$ws.Cells.Item(100, 13).Value = -1259.6666
$ws.Cells.Item(100, 14).Value = -4421.8
$ws.Cells.Item(132, 8).Value = 1829.3115
$ws.Cells.Item(132, 9).Value = 1786.4667
$ws.Cells.Item(132, 11).Value = 5359.4001
$ws.Cells.Item(132, 13).Value = -2829.4001
$ws.Cells.Item(135, 8).Value = 4171.0586
$ws.Cells.Item(135, 9).Value = 4734.0713
$ws.Cells.Item(135, 10).Value = 1543.6666
$ws.Cells.Item(135, 11).Value = 42606.64169999999
$ws.Cells.Item(135, 12).Value = 13892.9994
$ws.Cells.Item(135, 13).Value = -40071.64169999999
$ws.Cells.Item(135, 14).Value = -18962.9994
$ws.Cells.Item(137, 8).Value = 13359.3
$ws.Cells.Item(137, 9).Value = 5678
$ws.Cells.Item(137, 11).Value = 17034
$ws.Cells.Item(137, 13).Value = -14484

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 2571.2
$ws.Cells.Item(2, 9).Value = 1320
$ws.Cells.Item(2, 11).Value = 1320
$ws.Cells.Item(2, 13).Value = -1207
$ws.Cells.Item(32, 8).Value = 4830.0327
$ws.Cells.Item(32, 9).Value = 4905.68
$ws.Cells.Item(32, 11).Value = 4905.68
$ws.Cells.Item(32, 13).Value = -4618.68
$ws.Cells.Item(45, 8).Value = 1416.5
$ws.Cells.Item(45, 9).Value = 1409
$ws.Cells.Item(45, 11).Value = 1409
$ws.Cells.Item(45, 13).Value = -1032
$ws.Cells.Item(61, 8).Value = 3115.2778
$ws.Cells.Item(61, 9).Value = 2375.1155
$ws.Cells.Item(61, 11).Value = 2375.1155
$ws.Cells.Item(61, 13).Value = -2163.1155
$ws.Cells.Item(116, 8).Value = 2571.2
$ws.Cells.Item(116, 9).Value = 1320
$ws.Cells.Item(116, 11).Value = 1320
$ws.Cells.Item(116, 13).Value = 974
$ws.Cells.Item(122, 8).Value = 3154.5789
$ws.Cells.Item(122, 9).Value = 1381.0769
$ws.Cells.Item(122, 11).Value = 4143.2307
$ws.Cells.Item(122, 13).Value = -1693.2307
$ws.Cells.Item(132, 8).Value = 88683.5
$ws.Cells.Item(132, 9).Value = 2759.8
$ws.Cells.Item(132, 11).Value = 8279.400000000001
$ws.Cells.Item(132, 13).Value = -5749.400000000001
$ws.Cells.Item(136, 8).Value = 3115.2778
$ws.Cells.Item(136, 9).Value = 2375.1155
$ws.Cells.Item(136, 11).Value = 7125.3465
$ws.Cells.Item(136, 13).Value = -4575.3465

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 2571.2
$ws.Cells.Item(3, 9).Value = 1320
$ws.Cells.Item(3, 11).Value = 1320
$ws.Cells.Item(3, 13).Value = -1206
$ws.Cells.Item(26, 8).Value = 16480.572
$ws.Cells.Item(26, 9).Value = 16480.572
$ws.Cells.Item(26, 11).Value = 16480.572
$ws.Cells.Item(26, 13).Value = -16188.572
$ws.Cells.Item(94, 8).Value = 2109.4644
$ws.Cells.Item(94, 9).Value = 1700.6
$ws.Cells.Item(94, 11).Value = 1700.6
$ws.Cells.Item(94, 13).Value = -1249.6
$ws.Cells.Item(105, 8).Value = 9346.571

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 4751.0435
$ws.Cells.Item(31, 9).Value = 2806.4814
$ws.Cells.Item(31, 10).Value = 7514.3687
$ws.Cells.Item(31, 11).Value = 2806.4814
$ws.Cells.Item(31, 12).Value = 7514.3687
$ws.Cells.Item(31, 13).Value = -2511.4814
$ws.Cells.Item(31, 14).Value = -8104.3687
$ws.Cells.Item(34, 8).Value = 4751.0435
$ws.Cells.Item(34, 9).Value = 2806.4814
$ws.Cells.Item(34, 10).Value = 7514.3687
$ws.Cells.Item(34, 11).Value = 2806.4814
$ws.Cells.Item(34, 12).Value = 7514.3687
$ws.Cells.Item(34, 13).Value = -2604.4814
$ws.Cells.Item(34, 14).Value = -7918.3687
$ws.Cells.Item(132, 8).Value = 3287
$ws.Cells.Item(132, 9).Value = 3310.9167
$ws.Cells.Item(132, 11).Value = 9932.750100000001
$ws.Cells.Item(132, 13).Value = -7402.750100000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(64, 8).Value = 13579.333
$ws.Cells.Item(64, 10).Value = 14839.25
$ws.Cells.Item(64, 12).Value = 44517.75
$ws.Cells.Item(64, 14).Value = -45057.75
$ws.Cells.Item(67, 8).Value = 13579.333
$ws.Cells.Item(67, 10).Value = 14839.25
$ws.Cells.Item(67, 12).Value = 44517.75
$ws.Cells.Item(67, 14).Value = -46389.75
$ws.Cells.Item(121, 8).Value = 240.55556
$ws.Cells.Item(121, 9).Value = 104.57143
$ws.Cells.Item(121, 11).Value = 313.71429
$ws.Cells.Item(121, 13).Value = 996.28571
$ws.Cells.Item(129, 8).Value = 1472.5834
$ws.Cells.Item(129, 10).Value = 1703.6666
$ws.Cells.Item(129, 12).Value = 5110.9998
$ws.Cells.Item(129, 14).Value = -15110.9998

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(24, 8).Value = 21370.75
$ws.Cells.Item(24, 9).Value = 19710.857
$ws.Cells.Item(24, 10).Value = 32990
$ws.Cells.Item(24, 11).Value = 19710.857
$ws.Cells.Item(24, 12).Value = 32990
$ws.Cells.Item(24, 13).Value = -19537.857
$ws.Cells.Item(24, 14).Value = -33336
$ws.Cells.Item(102, 8).Value = 2256
$ws.Cells.Item(102, 9).Value = 2408
$ws.Cells.Item(102, 11).Value = 2408
$ws.Cells.Item(102, 13).Value = -786
$ws.Cells.Item(126, 8).Value = 93401.37
$ws.Cells.Item(126, 9).Value = 168736.17
$ws.Cells.Item(126, 11).Value = 506208.51
$ws.Cells.Item(126, 13).Value = -503738.51
$ws.Cells.Item(132, 8).Value = 1893.4
$ws.Cells.Item(132, 9).Value = 1991.75
$ws.Cells.Item(132, 11).Value = 5975.25
$ws.Cells.Item(132, 13).Value = -3445.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(40, 8).Value = 4402.25
$ws.Cells.Item(40, 9).Value = 3896.0667
$ws.Cells.Item(40, 11).Value = 3896.0667
$ws.Cells.Item(40, 13).Value = -3760.0667
$ws.Cells.Item(61, 8).Value = 3577.739
$ws.Cells.Item(61, 9).Value = 3382
$ws.Cells.Item(61, 10).Value = 4132.3335
$ws.Cells.Item(61, 11).Value = 3382
$ws.Cells.Item(61, 12).Value = 4132.3335
$ws.Cells.Item(61, 13).Value = -3180
$ws.Cells.Item(61, 14).Value = -4536.3335
$ws.Cells.Item(113, 8).Value = 3577.739
$ws.Cells.Item(113, 9).Value = 3382
$ws.Cells.Item(113, 10).Value = 4132.3335
$ws.Cells.Item(113, 11).Value = 3382
$ws.Cells.Item(113, 12).Value = 4132.3335
$ws.Cells.Item(113, 13).Value = -1212
$ws.Cells.Item(113, 14).Value = -8472.333500000001
$ws.Cells.Item(122, 8).Value = 5965.364
$ws.Cells.Item(122, 9).Value = 3620.3333
$ws.Cells.Item(122, 11).Value = 10860.9999
$ws.Cells.Item(122, 13).Value = -8410.999899999999
$ws.Cells.Item(132, 8).Value = 4556.048
$ws.Cells.Item(132, 9).Value = 2428.8918
$ws.Cells.Item(132, 11).Value = 7286.6754
$ws.Cells.Item(132, 13).Value = -4756.6754

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(122, 8).Value = 51820.75
$ws.Cells.Item(122, 9).Value = 112230.89
$ws.Cells.Item(122, 11).Value = 336692.67
$ws.Cells.Item(122, 13).Value = -334242.67
$ws.Cells.Item(136, 8).Value = 11114.611
$ws.Cells.Item(136, 9).Value = 10059.223
$ws.Cells.Item(136, 10).Value = 12170
$ws.Cells.Item(136, 11).Value = 30177.669
$ws.Cells.Item(136, 12).Value = 36510
$ws.Cells.Item(136, 13).Value = -27627.669
$ws.Cells.Item(136, 14).Value = -41610
